$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.716.80"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.539.16"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9981"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3949"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3204"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07205"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.084"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.754"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.658"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "1.540.66"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9983"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.161"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.372"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "21.720.00"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.407"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.854"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "1.715.40"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.126"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9804"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.28%  "
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.594"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.220"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02244"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.488"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06007"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2054"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.186"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9977"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5848"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5610"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.170"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06750"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.94%  "
